$wb = $excel.ActiveWorkbook

# ---- Sheet "snapshot": update scraped_at (K) for existing unchanged rows 2-24 ----
$snap = $wb.Worksheets.Item("snapshot")

$snap.Range("K2").Value = "2025-12-09T03:01:35.746350+00:00"
$snap.Range("K3").Value = "2025-12-09T03:01:35.746384+00:00"
$snap.Range("K4").Value = "2025-12-09T03:01:35.746404+00:00"
$snap.Range("K5").Value = "2025-12-09T03:01:38.083302+00:00"
$snap.Range("K6").Value = "2025-12-09T03:01:38.083335+00:00"
$snap.Range("K7").Value = "2025-12-09T03:01:40.830147+00:00"
$snap.Range("K8").Value = "2025-12-09T03:01:43.625061+00:00"
$snap.Range("K9").Value = "2025-12-09T03:01:46.021544+00:00"
$snap.Range("K10").Value = "2025-12-09T03:01:48.782450+00:00"
$snap.Range("K11").Value = "2025-12-09T03:01:53.897838+00:00"
$snap.Range("K12").Value = "2025-12-09T03:01:53.897866+00:00"
$snap.Range("K13").Value = "2025-12-09T03:01:56.621346+00:00"
$snap.Range("K14").Value = "2025-12-09T03:01:59.325241+00:00"
$snap.Range("K15").Value = "2025-12-09T03:02:02.162303+00:00"
$snap.Range("K16").Value = "2025-12-09T03:02:04.884909+00:00"
$snap.Range("K17").Value = "2025-12-09T03:02:04.884942+00:00"
$snap.Range("K18").Value = "2025-12-09T03:02:07.698589+00:00"
$snap.Range("K19").Value = "2025-12-09T03:02:07.698618+00:00"
$snap.Range("K20").Value = "2025-12-09T03:02:07.698636+00:00"
$snap.Range("K21").Value = "2025-12-09T03:02:10.418862+00:00"
$snap.Range("K22").Value = "2025-12-09T03:02:10.418891+00:00"
$snap.Range("K23").Value = "2025-12-09T03:02:10.418910+00:00"
$snap.Range("K24").Value = "2025-12-09T03:02:10.418927+00:00"

# ---- Insert new row 25 (new injury: SKA / Sapego Sergei) and shift rows 25-37 down to 26-38 ----
$snap.Rows.Item(25).Insert()

$snap.Range("A25").Value = "СКА"
$snap.Range("B25").Value = "СКА"
$snap.Range("C25").Value = "ska"
$snap.Range("D25").Value = "Сапего Сергей"
$snap.Range("E25").Value = "'9"
$snap.Range("F25").Value = "защитник"
$snap.Range("G25").Value = "'39875"
$snap.Range("H25").Value = "1369_СКА_сапегосергей"
$snap.Range("I25").Value = "injured_active"
$snap.Range("J25").Value = "https://www.khl.ru/clubs/ska/team/"
$snap.Range("K25").Value = "2025-12-09T03:02:10.418944+00:00"

# row 26: Бикмуллин Рафаэль
$snap.Range("A26").Value = "СОЧ"
$snap.Range("B26").Value = "ХК Сочи"
$snap.Range("C26").Value = "hc_sochi"
$snap.Range("D26").Value = "Бикмуллин Рафаэль"
$snap.Range("E26").Value = "'24"
$snap.Range("F26").Value = "нападающий"
$snap.Range("G26").Value = "'22424"
$snap.Range("H26").Value = "1369_СОЧ_бикмуллинрафаэль"
$snap.Range("I26").Value = "injured_active"
$snap.Range("J26").Value = "https://www.khl.ru/clubs/hc_sochi/team/"
$snap.Range("K26").Value = "2025-12-09T03:02:13.160715+00:00"

# row 27: Венгрыжановский Денис
$snap.Range("A27").Value = "СОЧ"
$snap.Range("B27").Value = "ХК Сочи"
$snap.Range("C27").Value = "hc_sochi"
$snap.Range("D27").Value = "Венгрыжановский Денис"
$snap.Range("E27").Value = "'9"
$snap.Range("F27").Value = "нападающий"
$snap.Range("G27").Value = "'31892"
$snap.Range("H27").Value = "1369_СОЧ_венгрыжановскийденис"
$snap.Range("I27").Value = "injured_active"
$snap.Range("J27").Value = "https://www.khl.ru/clubs/hc_sochi/team/"
$snap.Range("K27").Value = "2025-12-09T03:02:13.160746+00:00"

# row 28: Хёфенмайер Ноэль
$snap.Range("A28").Value = "СОЧ"
$snap.Range("B28").Value = "ХК Сочи"
$snap.Range("C28").Value = "hc_sochi"
$snap.Range("D28").Value = "Хёфенмайер Ноэль"
$snap.Range("E28").Value = "'22"
$snap.Range("F28").Value = "защитник"
$snap.Range("G28").Value = "'44847"
$snap.Range("H28").Value = "1369_СОЧ_хефенмайерноэль"
$snap.Range("I28").Value = "injured_active"
$snap.Range("J28").Value = "https://www.khl.ru/clubs/hc_sochi/team/"
$snap.Range("K28").Value = "2025-12-09T03:02:13.160766+00:00"

# row 29: Вишневский Дмитрий
$snap.Range("A29").Value = "СПР"
$snap.Range("B29").Value = "Спартак"
$snap.Range("C29").Value = "spartak"
$snap.Range("D29").Value = "Вишневский Дмитрий"
$snap.Range("E29").Value = "'55"
$snap.Range("F29").Value = "защитник"
$snap.Range("G29").Value = "'15299"
$snap.Range("H29").Value = "1369_СПР_вишневскийдмитрий"
$snap.Range("I29").Value = "injured_active"
$snap.Range("J29").Value = "https://www.khl.ru/clubs/spartak/team/"
$snap.Range("K29").Value = "2025-12-09T03:02:15.915796+00:00"

# row 30: Воронин Кирилл
$snap.Range("A30").Value = "ТОР"
$snap.Range("B30").Value = "Торпедо"
$snap.Range("C30").Value = "torpedo"
$snap.Range("D30").Value = "Воронин Кирилл"
$snap.Range("E30").Value = "'41"
$snap.Range("F30").Value = "нападающий"
$snap.Range("G30").Value = "'17354"
$snap.Range("H30").Value = "1369_ТОР_воронинкирилл"
$snap.Range("I30").Value = "injured_active"
$snap.Range("J30").Value = "https://www.khl.ru/clubs/torpedo/team/"
$snap.Range("K30").Value = "2025-12-09T03:02:20.893500+00:00"

# row 31: Кручинин Алексей
$snap.Range("A31").Value = "ТОР"
$snap.Range("B31").Value = "Торпедо"
$snap.Range("C31").Value = "torpedo"
$snap.Range("D31").Value = "Кручинин Алексей"
$snap.Range("E31").Value = "'78"
$snap.Range("F31").Value = "нападающий"
$snap.Range("G31").Value = "'16355"
$snap.Range("H31").Value = "1369_ТОР_кручининалексей"
$snap.Range("I31").Value = "injured_active"
$snap.Range("J31").Value = "https://www.khl.ru/clubs/torpedo/team/"
$snap.Range("K31").Value = "2025-12-09T03:02:20.893527+00:00"

# row 32: Принс Шэйн
$snap.Range("A32").Value = "ТОР"
$snap.Range("B32").Value = "Торпедо"
$snap.Range("C32").Value = "torpedo"
$snap.Range("D32").Value = "Принс Шэйн"
$snap.Range("E32").Value = "'18"
$snap.Range("F32").Value = "нападающий"
$snap.Range("G32").Value = "'19060"
$snap.Range("H32").Value = "1369_ТОР_принсшэйн"
$snap.Range("I32").Value = "injured_active"
$snap.Range("J32").Value = "https://www.khl.ru/clubs/torpedo/team/"
$snap.Range("K32").Value = "2025-12-09T03:02:20.893544+00:00"

# row 33: Мыльников Сергей И
$snap.Range("A33").Value = "ТРК"
$snap.Range("B33").Value = "Трактор"
$snap.Range("C33").Value = "traktor"
$snap.Range("D33").Value = "Мыльников Сергей И"
$snap.Range("E33").Value = "'20"
$snap.Range("F33").Value = "вратарь"
$snap.Range("G33").Value = "'24799"
$snap.Range("H33").Value = "1369_ТРК_мыльниковсергейи"
$snap.Range("I33").Value = "injured_active"
$snap.Range("J33").Value = "https://www.khl.ru/clubs/traktor/team/"
$snap.Range("K33").Value = "2025-12-09T03:02:23.209992+00:00"

# row 34: Светлаков Андрей
$snap.Range("A34").Value = "ТРК"
$snap.Range("B34").Value = "Трактор"
$snap.Range("C34").Value = "traktor"
$snap.Range("D34").Value = "Светлаков Андрей"
$snap.Range("E34").Value = "'87"
$snap.Range("F34").Value = "нападающий"
$snap.Range("G34").Value = "'19218"
$snap.Range("H34").Value = "1369_ТРК_светлаковандрей"
$snap.Range("I34").Value = "injured_active"
$snap.Range("J34").Value = "https://www.khl.ru/clubs/traktor/team/"
$snap.Range("K34").Value = "2025-12-09T03:02:23.210021+00:00"

# row 35: Бучельников Дмитрий
$snap.Range("A35").Value = "ЦСК"
$snap.Range("B35").Value = "ЦСКА"
$snap.Range("C35").Value = "cska"
$snap.Range("D35").Value = "Бучельников Дмитрий"
$snap.Range("E35").Value = "'72"
$snap.Range("F35").Value = "нападающий"
$snap.Range("G35").Value = "'39102"
$snap.Range("H35").Value = "1369_ЦСК_бучельниковдмитрий"
$snap.Range("I35").Value = "injured_active"
$snap.Range("J35").Value = "https://www.khl.ru/clubs/cska/team/"
$snap.Range("K35").Value = "2025-12-09T03:02:25.598156+00:00"

# row 36: Моисеев Данила
$snap.Range("A36").Value = "ЦСК"
$snap.Range("B36").Value = "ЦСКА"
$snap.Range("C36").Value = "cska"
$snap.Range("D36").Value = "Моисеев Данила"
$snap.Range("E36").Value = "'93"
$snap.Range("F36").Value = "нападающий"
$snap.Range("G36").Value = "'23931"
$snap.Range("H36").Value = "1369_ЦСК_моисеевданила"
$snap.Range("I36").Value = "injured_active"
$snap.Range("J36").Value = "https://www.khl.ru/clubs/cska/team/"
$snap.Range("K36").Value = "2025-12-09T03:02:25.598184+00:00"

# row 37: Бишофф Джейк
$snap.Range("A37").Value = "ШДР"
$snap.Range("B37").Value = "Драконы"
$snap.Range("C37").Value = "kunlun"
$snap.Range("D37").Value = "Бишофф Джейк"
$snap.Range("E37").Value = "'28"
$snap.Range("F37").Value = "защитник"
$snap.Range("G37").Value = "'45490"
$snap.Range("H37").Value = "1369_ШДР_бишоффджейк"
$snap.Range("I37").Value = "injured_active"
$snap.Range("J37").Value = "https://www.khl.ru/clubs/kunlun/team/"
$snap.Range("K37").Value = "2025-12-09T03:02:28.382069+00:00"

# row 38: Гроло Жереми
$snap.Range("A38").Value = "ШДР"
$snap.Range("B38").Value = "Драконы"
$snap.Range("C38").Value = "kunlun"
$snap.Range("D38").Value = "Гроло Жереми"
$snap.Range("E38").Value = "'75"
$snap.Range("F38").Value = "защитник"
$snap.Range("G38").Value = "'45343"
$snap.Range("H38").Value = "1369_ШДР_гроложереми"
$snap.Range("I38").Value = "injured_active"
$snap.Range("J38").Value = "https://www.khl.ru/clubs/kunlun/team/"
$snap.Range("K38").Value = "2025-12-09T03:02:28.382096+00:00"

# ---- Sheet "new_injured": remove the old Murphy Trevor row, update Sapego Sergei row ----
$newinj = $wb.Worksheets.Item("new_injured")
$newinj.Rows.Item(3).Delete()
$newinj.Range("C2").Value = "Сапего Сергей"
$newinj.Range("D2").Value = "1369_СКА_сапегосергей"
$newinj.Range("F2").Value = "2025-12-09T11:02:28.891776+08:00"
$newinj.Range("G2").Value = "'2025-12-09"
